$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the pre-existing "_GoBack" bookmark (it sat right after the
#    class-diagram image, marking where the previous author last edited).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldGoBack = $d.Bookmarks.Item("_GoBack")
    $oldGoBack.Delete()
}

# ---------------------------------------------------------------------
# 2. Replace the word "workshop" with "milestone" in the submission-
#    instructions paragraph, keeping the trailing period untouched.
# ---------------------------------------------------------------------
$wordRange = $d.Content
$found = $wordRange.Find.Execute("workshop", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the word 'workshop' to replace."
}
$milestoneStart = $wordRange.Start
$replacement = "milestone"
$wordRange.Text = $replacement
$milestoneEnd = $milestoneStart + $replacement.Length

# ---------------------------------------------------------------------
# 3. Recreate the run boundary that sits right before "milestone" (Word
#    always breaks a run at a bookmark), then drop a fresh "_GoBack"
#    bookmark at the point right after "milestone" -- i.e. where the
#    cursor was left after typing the replacement word -- and finally
#    discard the temporary helper bookmark used only to force the left
#    run boundary.
# ---------------------------------------------------------------------
$splitPoint = $d.Range($milestoneStart, $milestoneStart)
$d.Bookmarks.Add("__TempRunSplit", $splitPoint)

$goBackPoint = $d.Range($milestoneEnd, $milestoneEnd)
$d.Bookmarks.Add("_GoBack", $goBackPoint)

$tempBm = $d.Bookmarks.Item("__TempRunSplit")
$tempBm.Delete()
